$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sheet1")

# Insert a new column at BO, shifting the existing BO column (the numeric
# "order" value used by cutsendmail) one position to the right, into BP.
# Using Range/Columns Insert (rather than building cells from scratch)
# keeps the original cell formatting/style on both the old and new cells.
$ws.Columns("BO").Insert(-4161)   # -4161 = xlShiftToRight

# The freed-up BO column becomes the new "group" column (the mailing-list
# name that will be used to build the json file for cutsendmail). Its
# value is simply a copy of the existing group column BN for each data row.
for ($r = 2; $r -le 15; $r++) {
    $group = $ws.Cells.Item($r, 66).Value2   # column BN = 66 -> group name
    if ($group -ne $null -and $group -ne "") {
        $ws.Cells.Item($r, 67).Value = $group  # column BO = 67
    }
}
